$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the percentage values in column F (rows 19, 21, 22, 23)
$ws.Range("F19").Value = 0.15
$ws.Range("F21").Value = 0.15
$ws.Range("F22").Value = 0.2
$ws.Range("F23").Value = 0.2

# Add a new total row with a SUM formula
$ws.Range("F25").Formula = "=SUM(F16:F23)"

# Update the view: scroll so A10 is the top-left visible cell, and
# change the active selection to H27
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 10
$ws.Range("H27").Select()
